# Change ways of setting gradient angle:
# replace the fully-transparent "rotated linear gradient" fills
# (gradFill rotWithShape, two black/alpha-0 stops, ang=16200000) with an
# equivalent fully-transparent solidFill using the shape's accent color,
# and recolor the matching text runs from black to the same accent color.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# RGB() equivalent for PowerPoint long color values (0xBBGGRR)
# 00BCA7 -> teal accent ; 00584E -> dark teal accent
$teal = 10992640   # 0x00BCA7
$darkTeal = 5134336  # 0x00584E

function Set-SolidAlpha0 {
    param($shape, $rgb)
    $shape.Fill.Solid()
    $shape.Fill.Transparency = 1
    $shape.Fill.ForeColor.RGB = $rgb
    $shape.Fill.Transparency = 1
}

# --- Shapes whose gradFill becomes a transparent solidFill -----------------
# 1  id=44 "任意多边形: 形状 43" -> teal
Set-SolidAlpha0 $s.Shapes.Item(1) $teal
# 7  id=38 "矩形 37" (full-bleed backdrop rect) -> dark teal
Set-SolidAlpha0 $s.Shapes.Item(7) $darkTeal
# 8  id=9  "rectangle-11" (badge behind 关键词1) -> teal
Set-SolidAlpha0 $s.Shapes.Item(8) $teal
# 13 id=15 "椭圆 14" -> dark teal
Set-SolidAlpha0 $s.Shapes.Item(13) $darkTeal
# 14 id=64 "椭圆 63" -> teal
Set-SolidAlpha0 $s.Shapes.Item(14) $teal
# 15 id=4  "rectangle-11" (badge behind 关键词2) -> teal
Set-SolidAlpha0 $s.Shapes.Item(15) $teal
# 17 id=7  "rectangle-11" (badge behind 关键词3) -> teal
Set-SolidAlpha0 $s.Shapes.Item(17) $teal
# 19 id=10 "rectangle-11" (badge behind 关键词4) -> teal
Set-SolidAlpha0 $s.Shapes.Item(19) $teal
# 21 id=18 "rectangle-11" (badge behind 关键词5) -> teal
Set-SolidAlpha0 $s.Shapes.Item(21) $teal
# 23 id=20 "rectangle-11" (badge behind 关键词6) -> teal
Set-SolidAlpha0 $s.Shapes.Item(23) $teal
# 25 id=22 "rectangle-11" (badge behind 关键词7) -> teal
Set-SolidAlpha0 $s.Shapes.Item(25) $teal
# 27 id=24 "rectangle-11" (badge behind 关键词8) -> teal
Set-SolidAlpha0 $s.Shapes.Item(27) $teal

# --- Keyword label runs: black -> dark teal ---------------------------------
$s.Shapes.Item(9).TextFrame.TextRange.Font.Color.RGB = $darkTeal   # 关键词1
$s.Shapes.Item(16).TextFrame.TextRange.Font.Color.RGB = $darkTeal  # 关键词2
$s.Shapes.Item(18).TextFrame.TextRange.Font.Color.RGB = $darkTeal  # 关键词3
$s.Shapes.Item(20).TextFrame.TextRange.Font.Color.RGB = $darkTeal  # 关键词4
$s.Shapes.Item(22).TextFrame.TextRange.Font.Color.RGB = $darkTeal  # 关键词5
$s.Shapes.Item(24).TextFrame.TextRange.Font.Color.RGB = $darkTeal  # 关键词6
$s.Shapes.Item(26).TextFrame.TextRange.Font.Color.RGB = $darkTeal  # 关键词7
$s.Shapes.Item(28).TextFrame.TextRange.Font.Color.RGB = $darkTeal  # 关键词8

# --- Title runs: black -> teal ----------------------------------------------
$s.Shapes.Item(29).TextFrame.TextRange.Font.Color.RGB = $teal      # 标题

# --- Description runs: black -> dark teal -----------------------------------
$s.Shapes.Item(30).TextFrame.TextRange.Font.Color.RGB = $darkTeal  # 关键词1-简述
$s.Shapes.Item(31).TextFrame.TextRange.Font.Color.RGB = $darkTeal  # 关键词2-简述
$s.Shapes.Item(32).TextFrame.TextRange.Font.Color.RGB = $darkTeal  # 关键词3-简述
$s.Shapes.Item(33).TextFrame.TextRange.Font.Color.RGB = $darkTeal  # 关键词4-简述
$s.Shapes.Item(34).TextFrame.TextRange.Font.Color.RGB = $darkTeal  # 关键词5-简述
$s.Shapes.Item(35).TextFrame.TextRange.Font.Color.RGB = $darkTeal  # 关键词6-简述
$s.Shapes.Item(36).TextFrame.TextRange.Font.Color.RGB = $darkTeal  # 关键词7-简述
$s.Shapes.Item(37).TextFrame.TextRange.Font.Color.RGB = $darkTeal  # 关键词8-简述
